$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# <w:tblInd w:w="81"/> -> <w:tblInd w:w="72"/>  (81/20=4.05pt -> 72/20=3.6pt)
$t.Rows.LeftIndent = 3.6

# <w:tblCellMar><w:left w:w="70"/></w:tblCellMar> -> w:w="60" (3.5pt -> 3.0pt)
$t.LeftPadding = 3.0

# first column: <w:gridCol w:w="4014"/> / <w:tcW w:w="4014"/> -> 4013 (200.7pt -> 200.65pt)
$t.Columns(1).Width = 200.65

# every cell's <w:tcMar><w:left w:w="70"/></w:tcMar> -> w:w="60" (3.5pt -> 3.0pt)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.LeftPadding = 3.0
    }
}

# "maxim 10000 caractere" -> "maxim 1000 caractere" in the constraints paragraph
$d.Content.Find.Execute("Fiecare ecuaţie e formată din maxim 10000 caractere;", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fiecare ecuaţie e formată din maxim 1000 caractere;", 2)
